# GNTW Yearly Financials update: add a new latest fiscal-year column.
# A brand new data column is inserted immediately to the left of the
# existing "D" column (the most-recent-year column), pushing the
# historical D:K columns one slot to the right (E:L), and the new D
# column is populated with the newly reported FY2018 figures for every
# financial statement line across the three statements
# (Income Statement / Balance Sheet / Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before D; existing D:K data slides to E:L.
$ws.Columns("D").Insert()

# The freshly inserted column has no number formatting / style of its
# own yet - copy it over from column E (which holds what used to be
# column D, i.e. the adjacent, correctly-formatted column) so the new
# column matches (date format for the header row, number format for
# the data rows, etc.)
$ws.Columns("E").Copy()
$ws.Columns("D").PasteSpecial(-4122)

# --- Income Statement -------------------------------------------------
$ws.Range("D7").Value = 43343
$ws.Range("D8").Value = "NA"
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("D18").Value = "NA"
$ws.Range("D20").Value = "NA"
$ws.Range("D21").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = "NA"
$ws.Range("D33").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 0

# --- Balance Sheet ------------------------------------------------------
$ws.Range("D38").Value = 43343
$ws.Range("D41").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 0
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 100
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 200
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -200
$ws.Range("D77").Value = 0

# --- Cash Flow Statement -------------------------------------------------
$ws.Range("D80").Value = 43343
$ws.Range("D81").Value = 0
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0
